# Apply updated crypto price/volume snapshot values (GitHub Actions refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    # Force the cell to stay a text cell (matches the source data's inlineStr/
    # shared-string type) even when the value looks like a number (e.g. "41.02"),
    # then restore the default "Normal" style so no stray number format sticks.
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

# Row 2
Set-TextValue $ws.Range("D2") '34.690.19'
$ws.Range("E2").Value = '  -2.41%  '

# Row 3
Set-TextValue $ws.Range("D3") '1.876.59'
$ws.Range("E3").Value = '  -1.91%  '

# Row 4
$ws.Range("E4").Value = '  -0.77%  '

# Row 5
Set-TextValue $ws.Range("D5") '247.72'
$ws.Range("E5").Value = '  +0.34%  '

# Row 6
Set-TextValue $ws.Range("D6") '0.687'
$ws.Range("E6").Value = '  -2.87%  '

# Row 7
$ws.Range("E7").Value = '  -0.91%  '

# Row 8
Set-TextValue $ws.Range("D8") '41.02'
$ws.Range("E8").Value = '  +0.29%  '

# Row 9
Set-TextValue $ws.Range("D9") '0.346'
$ws.Range("E9").Value = '  -2.64%  '

# Row 10
Set-TextValue $ws.Range("D10") '50.78'
$ws.Range("E10").Value = '  -3.57%  '

# Row 11
Set-TextValue $ws.Range("D11") '0.0736'
$ws.Range("E11").Value = '  +0.58%  '

# Row 12
Set-TextValue $ws.Range("D12") '0.0968'
$ws.Range("E12").Value = '  -2.15%  '

# Row 13
Set-TextValue $ws.Range("D13") '2.150.71'
$ws.Range("E13").Value = '  -1.83%  '

# Row 14
Set-TextValue $ws.Range("D14") '12.78'
$ws.Range("E14").Value = '  +1.24%  '

# Row 15
Set-TextValue $ws.Range("D15") '0.713'
$ws.Range("E15").Value = '  -0.49%  '

# Row 16
Set-TextValue $ws.Range("D16") '4.88'
$ws.Range("E16").Value = '  -0.74%  '

# Row 17
Set-TextValue $ws.Range("D17") '1.882.61'
$ws.Range("E17").Value = '  -1.50%  '

# Row 18
Set-TextValue $ws.Range("D18") '34.723.27'
$ws.Range("E18").Value = '  -2.27%  '

# Row 19
Set-TextValue $ws.Range("D19") '72.76'
$ws.Range("E19").Value = '  -0.74%  '

# Row 20
Set-TextValue $ws.Range("D20") '0.0₃0820'
$ws.Range("E20").Value = '  -0.98%  '

# Row 21
Set-TextValue $ws.Range("D21") '247.56'
$ws.Range("E21").Value = '  +1.95%  '

# Row 22
Set-TextValue $ws.Range("D22") '12.71'
$ws.Range("E22").Value = '  -3.42%  '

# Row 23
Set-TextValue $ws.Range("D23") '4.90'
$ws.Range("E23").Value = '  -2.74%  '

# Row 24
$ws.Range("E24").Value = '  -0.86%  '

# Row 25
$ws.Range("E25").Value = '  +2.36%  '

# Row 26
$ws.Range("E26").Value = '  -4.14%  '

# Row 27
Set-TextValue $ws.Range("D27") '164.86'
$ws.Range("E27").Value = '  -2.65%  '

# Row 28
Set-TextValue $ws.Range("D28") '8.36'
$ws.Range("E28").Value = '  -3.52%  '

# Row 29
Set-TextValue $ws.Range("D29") '18.19'
$ws.Range("E29").Value = '  -3.19%  '

# Row 30
Set-TextValue $ws.Range("D30") '0.127'
$ws.Range("E30").Value = '  -3.96%  '

# Row 31
Set-TextValue $ws.Range("D31") '4.128.63'
$ws.Range("E31").Value = '  +0.21%  '

# Row 32
$ws.Range("E32").Value = '  +0.18%  '

# Row 33
Set-TextValue $ws.Range("D33") '0.0579'
$ws.Range("E33").Value = '  +0.62%  '

# Row 34
Set-TextValue $ws.Range("D34") '1.54'
$ws.Range("E34").Value = '  +3.51%  '

# Row 35
Set-TextValue $ws.Range("D35") '4.14'
$ws.Range("E35").Value = '  -1.89%  '

# Row 36
$ws.Range("E36").Value = '  -0.83%  '

# Row 37
Set-TextValue $ws.Range("D37") '1.83'
$ws.Range("E37").Value = '  -2.20%  '

# Row 38
Set-TextValue $ws.Range("D38") '0.829'
$ws.Range("E38").Value = '  -9.52%  '

# Row 39
$ws.Range("E39").Value = '  -3.17%  '

# Row 40
Set-TextValue $ws.Range("D40") '17.28'
$ws.Range("E40").Value = '  +0.36%  '

# Row 41
Set-TextValue $ws.Range("D41") '97.71'
$ws.Range("E41").Value = '  -0.06%  '

# Row 42
$ws.Range("E42").Value = '  +1.76%  '

# Row 43
$ws.Range("E43").Value = '  -0.23%  '

# Row 44
$ws.Range("E44").Value = '  -3.89%  '

# Row 45
Set-TextValue $ws.Range("D45") '1.289.72'
$ws.Range("E45").Value = '  -4.93%  '

# Row 46
Set-TextValue $ws.Range("D46") '2.35'
$ws.Range("E46").Value = '  -3.97%  '

# Row 47
$ws.Range("E47").Value = '  -0.85%  '

# Row 48
Set-TextValue $ws.Range("D48") '2.72'
$ws.Range("E48").Value = '  -2.49%  '

# Row 49
Set-TextValue $ws.Range("D49") '12.15'
$ws.Range("E49").Value = '  -0.54%  '

# Row 50
Set-TextValue $ws.Range("D50") '0.0759'
$ws.Range("E50").Value = '  +5.63%  '

# Row 51
Set-TextValue $ws.Range("D51") '6.44'
$ws.Range("E51").Value = '  -1.69%  '
